$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 24.916566
$ws.Range("H2").Value = 74.749698
$ws.Range("I2").Value = 0.459912889255076
$ws.Range("J2").Value = 0.459912889255076
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 4.847498666666667
$ws.Range("N2").Value = 14.542496
$ws.Range("O2").Value = 0.03400671694637637
$ws.Range("P2").Value = 0.03400671694637637
$ws.Range("Q2").Value = 120.783020462912
$ws.Range("R2").Value = 1087.047184166208
$ws.Range("S2").Value = 0.01564012744488751
$ws.Range("T2").Value = 0.01564012744488751

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 24.916566
$ws.Range("H3").Value = 74.749698
$ws.Range("I3").Value = 0.459912889255076
$ws.Range("J3").Value = 0.459912889255076
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 1.627093333333334
$ws.Range("N3").Value = 4.88128
$ws.Range("O3").Value = 0.01141456784970118
$ws.Range("P3").Value = 0.01141456784970118
$ws.Range("Q3").Value = 40.54157842816
$ws.Range("R3").Value = 364.87420585344
$ws.Range("S3").Value = 0.005249706879354171
$ws.Range("T3").Value = 0.005249706879354171

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 24.916566
$ws.Range("H4").Value = 74.749698
$ws.Range("I4").Value = 0.459912889255076
$ws.Range("J4").Value = 0.459912889255076
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 136.0707373333333
$ws.Range("N4").Value = 408.212212
$ws.Range("O4").Value = 0.9545787152039225
$ws.Range("P4").Value = 0.9545787152039225
$ws.Range("Q4").Value = 3390.415507434664
$ws.Range("R4").Value = 30513.73956691198
$ws.Range("S4").Value = 0.4390230549308343
$ws.Range("T4").Value = 0.4390230549308343

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 19.60300333333333
$ws.Range("H5").Value = 58.80901
$ws.Range("I5").Value = 0.3618345281251927
$ws.Range("J5").Value = 0.3618345281251927
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 4.847498666666667
$ws.Range("N5").Value = 14.542496
$ws.Range("O5").Value = 0.03400671694637637
$ws.Range("P5").Value = 0.03400671694637637
$ws.Range("Q5").Value = 95.02553252099555
$ws.Range("R5").Value = 855.22979268896
$ws.Range("S5").Value = 0.01230480437937909
$ws.Range("T5").Value = 0.01230480437937909

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 19.60300333333333
$ws.Range("H6").Value = 58.80901
$ws.Range("I6").Value = 0.3618345281251927
$ws.Range("J6").Value = 0.3618345281251927
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 1.627093333333334
$ws.Range("N6").Value = 4.88128
$ws.Range("O6").Value = 0.01141456784970118
$ws.Range("P6").Value = 0.01141456784970118
$ws.Range("Q6").Value = 31.89591603697778
$ws.Range("R6").Value = 287.0632443328
$ws.Range("S6").Value = 0.004130184771649622
$ws.Range("T6").Value = 0.004130184771649622

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 19.60300333333333
$ws.Range("H7").Value = 58.80901
$ws.Range("I7").Value = 0.3618345281251927
$ws.Range("J7").Value = 0.3618345281251927
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 136.0707373333333
$ws.Range("N7").Value = 408.212212
$ws.Range("O7").Value = 0.9545787152039225
$ws.Range("P7").Value = 0.9545787152039225
$ws.Range("Q7").Value = 2667.395117514458
$ws.Range("R7").Value = 24006.55605763012
$ws.Range("S7").Value = 0.345399538974164
$ws.Range("T7").Value = 0.345399538974164

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 9.657138000000002
$ws.Range("H8").Value = 28.971414
$ws.Range("I8").Value = 0.1782525826197313
$ws.Range("J8").Value = 0.1782525826197313
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 4.847498666666667
$ws.Range("N8").Value = 14.542496
$ws.Range("O8").Value = 0.03400671694637637
$ws.Range("P8").Value = 0.03400671694637637
$ws.Range("Q8").Value = 46.81296357881601
$ws.Range("R8").Value = 421.3166722093441
$ws.Range("S8").Value = 0.006061785122109769
$ws.Range("T8").Value = 0.006061785122109769

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 9.657138000000002
$ws.Range("H9").Value = 28.971414
$ws.Range("I9").Value = 0.1782525826197313
$ws.Range("J9").Value = 0.1782525826197313
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 1.627093333333334
$ws.Range("N9").Value = 4.88128
$ws.Range("O9").Value = 0.01141456784970118
$ws.Range("P9").Value = 0.01141456784970118
$ws.Range("Q9").Value = 15.71306485888
$ws.Range("R9").Value = 141.41758372992
$ws.Range("S9").Value = 0.002034676198697388
$ws.Range("T9").Value = 0.002034676198697388

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 9.657138000000002
$ws.Range("H10").Value = 28.971414
$ws.Range("I10").Value = 0.1782525826197313
$ws.Range("J10").Value = 0.1782525826197313
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 136.0707373333333
$ws.Range("N10").Value = 408.212212
$ws.Range("O10").Value = 0.9545787152039225
$ws.Range("P10").Value = 0.9545787152039225
$ws.Range("Q10").Value = 1314.053888189752
$ws.Range("R10").Value = 11826.48499370777
$ws.Range("S10").Value = 0.1701561212989241
$ws.Range("T10").Value = 0.1701561212989241
